$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$formula = "=IF(ISBLANK(Tableau1[[#This Row],[Heure Début]]),`"`",Tableau1[[#This Row],[Heure fin]]-Tableau1[[#This Row],[Heure Début]])"

# Materialize the 12 new rows (29-40) with the same "empty template" pattern
# already used by row 28, by seeding the H-column formula (so each row exists)
# and then copying the formatting of row 28 down across them.
for ($r = 29; $r -le 40; $r++) {
    $ws.Range("H$r").Formula = $formula
}
$ws.Range("E28:M28").Copy()
$ws.Range("E29:M40").PasteSpecial(-4122)

# Row 28: new journal entry - "Légendes"
# (shared-string insertion order: J28 "Légendes" then L28 "Ajout de légende...")
$ws.Range("E28").Value = 44266
$ws.Range("F28").Value = 0.74305555555555547
$ws.Range("G28").Value = 0.75694444444444453
$ws.Range("I28").Value = "Développement"
$ws.Range("J28").Value = "Légendes"
$ws.Range("K28").Value = "Domicile"
$ws.Range("L28").Value = "Ajout de légende sur les côté de la grille"

# Row 29: new journal entry - option to quit
# (shared-string insertion order: L29 "Ajout d'une option pour quitter..."
# then J29 "ajout d'une option de quitter")
$ws.Range("E29").Value = 44266
$ws.Range("F29").Value = 0.76388888888888884
$ws.Range("G29").Value = 0.77777777777777779
$ws.Range("I29").Value = "Développement"
$ws.Range("L29").Value = "Ajout d'une option pour quitter pendant le jeux"
$ws.Range("K29").Value = "Domicile"
$ws.Range("J29").Value = "ajout d'une option de quitter"

# Give E28:G29 back their date/time number formats + alignment (copied from
# row 27, which already carries the correct date/time styles), since the
# plain Value assignment above reset them to the generic template style.
$ws.Range("E27:G27").Copy()
$ws.Range("E28:G29").PasteSpecial(-4122)

# Re-apply the values (PasteSpecial formats doesn't touch them, but make sure
# nothing drifted)
$ws.Range("E28").Value = 44266
$ws.Range("F28").Value = 0.74305555555555547
$ws.Range("G28").Value = 0.75694444444444453
$ws.Range("E29").Value = 44266
$ws.Range("F29").Value = 0.76388888888888884
$ws.Range("G29").Value = 0.77777777777777779

# Grow the table (Tableau1) to cover the new rows
$tbl = $ws.ListObjects.Item("Tableau1")
$tbl.Resize($ws.Range("E5:M40"))

$ws.Range("J30").Select()
